# Updates crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "35.259.78"
Set-TextValue "E2" "  -0.29%  "

Set-TextValue "D3" "1.909.22"
Set-TextValue "E3" "  +0.19%  "

Set-TextValue "E4" "  -0.07%  "

Set-TextValue "D5" "0.723"
Set-TextValue "E5" "  +9.16%  "

Set-TextValue "D6" "256.25"
Set-TextValue "E6" "  +4.02%  "

Set-TextValue "E7" "  -0.01%  "

Set-TextValue "D8" "40.61"
Set-TextValue "E8" "  -2.18%  "

Set-TextValue "D9" "0.375"
Set-TextValue "E9" "  +7.76%  "

Set-TextValue "D10" "52.90"
Set-TextValue "E10" "  -0.20%  "

Set-TextValue "D11" "0.0761"
Set-TextValue "E11" "  +5.39%  "

Set-TextValue "E12" "  -0.37%  "

Set-TextValue "D13" "2.186.35"
Set-TextValue "E13" "  +0.16%  "

Set-TextValue "D14" "12.85"
Set-TextValue "E14" "  +6.53%  "

Set-TextValue "D15" "0.727"
Set-TextValue "E15" "  +4.17%  "

Set-TextValue "B16" "Polkadot"
Set-TextValue "C16" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D16" "4.97"
Set-TextValue "E16" "  +2.21%  "

Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "1.903.62"
Set-TextValue "E17" "  -0.18%  "

Set-TextValue "D18" "35.240.84"
Set-TextValue "E18" "  -0.36%  "

Set-TextValue "D19" "74.91"
Set-TextValue "E19" "  +3.75%  "

Set-TextValue "E20" "  +3.39%  "

Set-TextValue "D21" "243.77"
Set-TextValue "E21" "  +1.38%  "

Set-TextValue "D22" "13.06"
Set-TextValue "E22" "  +4.58%  "

Set-TextValue "E23" "  +5.55%  "

Set-TextValue "E24" "  +0.00%  "

Set-TextValue "D25" "2.45"
Set-TextValue "E25" "  +7.11%  "

Set-TextValue "D26" "2.43"
Set-TextValue "E26" "  +3.37%  "

Set-TextValue "D27" "166.23"
Set-TextValue "E27" "  -2.44%  "

Set-TextValue "E28" "  +3.17%  "

Set-TextValue "D29" "18.74"
Set-TextValue "E29" "  +1.83%  "

Set-TextValue "D30" "0.133"
Set-TextValue "E30" "  +4.08%  "

Set-TextValue "D31" "4.128.96"
Set-TextValue "E31" "  +19.46%  "

Set-TextValue "D32" "4.39"
Set-TextValue "E32" "  +6.08%  "

Set-TextValue "E33" "  +14.77%  "

Set-TextValue "E34" "  +21.94%  "

Set-TextValue "D35" "0.0588"
Set-TextValue "E35" "  +3.89%  "

Set-TextValue "D36" "4.24"
Set-TextValue "E36" "  +3.33%  "

Set-TextValue "E37" "  -0.83%  "

Set-TextValue "E38" "  -1.77%  "

Set-TextValue "E39" "  -0.18%  "

Set-TextValue "B40" "InjectiveProtocol"
Set-TextValue "C40" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D40" "17.20"
Set-TextValue "E40" "  +5.14%  "

Set-TextValue "B41" "VeChain"
Set-TextValue "C41" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D41" "0.0218"
Set-TextValue "E41" "  +4.92%  "

Set-TextValue "D42" "96.42"
Set-TextValue "E42" "  +7.20%  "

Set-TextValue "E43" "  +1.20%  "

Set-TextValue "E44" "  +3.95%  "

Set-TextValue "D45" "1.336.29"
Set-TextValue "E45" "  -0.23%  "

Set-TextValue "D46" "2.43"
Set-TextValue "E46" "  +1.48%  "

Set-TextValue "D47" "2.43"
Set-TextValue "E47" "  +0.94%  "

Set-TextValue "E49" "  -0.49%  "

Set-TextValue "D50" "44.96"
Set-TextValue "E50" "  -6.86%  "

Set-TextValue "D51" "0.0753"
Set-TextValue "E51" "  +6.76%  "
